# Fruta / hortaliza, semanal
# Insert a new weekly record at the top of the data (row 9), pushing the
# existing history rows (old rows 9-16) down to rows 10-17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 9:16 down to 10:17, leaving a blank (format-inherited) row 9.
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with this week's Chirimoya market record.
$ws.Cells.Item(9, 1).Value = 10
$ws.Cells.Item(9, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(9, 3).Value = 'La Araucanía'
$ws.Cells.Item(9, 4).Value = 44413
$ws.Cells.Item(9, 5).Value = 9
$ws.Cells.Item(9, 6).Value = 'Fruta'
$ws.Cells.Item(9, 7).Value = 100107
$ws.Cells.Item(9, 8).Value = 'Otros'
$ws.Cells.Item(9, 9).Value = 100107002
$ws.Cells.Item(9, 10).Value = 'Chirimoya'
$ws.Cells.Item(9, 11).Value = 'Cultivar IV Región'
$ws.Cells.Item(9, 12).Value = 'Primera'
$ws.Cells.Item(9, 13).Value = 35
$ws.Cells.Item(9, 14).Value = 3500
$ws.Cells.Item(9, 15).Value = 3500
$ws.Cells.Item(9, 16).Value = 3500
$ws.Cells.Item(9, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(9, 18).Value = 'Provincia del Elquí'
$ws.Cells.Item(9, 19).Value = 3500
$ws.Cells.Item(9, 20).Value = 1

Write-Host "New dimension:" $ws.UsedRange.Address()
